$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.007.41"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "1.636.80"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.57"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.85"
$ws.Range("E8").Value = "  +11.09%  "
$ws.Range("E9").Value = "  +4.22%  "
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.870.56"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("D13").Value = "1.642.68"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.71"
$ws.Range("E14").Value = "  +27.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.580"
$ws.Range("E15").Value = "  +7.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.92"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").Value = "30.022.85"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.89"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.73"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").Value = "0.0₃0710"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +5.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.72"
$ws.Range("E23").Value = "  +5.15%  "
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.31"
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.73"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0491"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("E31").Value = "  +6.48%  "
$ws.Range("E32").Value = "  +5.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").Value = "1.432.66"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  +7.84%  "
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.80"
$ws.Range("E41").Value = "  +12.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.839"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.14"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0499"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").Value = "1.776.86"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.49"
$ws.Range("E50").Value = "  +5.01%  "
$ws.Range("D51").Value = "0.0₆0110"
$ws.Range("E51").Value = "  +5.32%  "
